# Insert a new data row at row 85 (pushing existing rows 85..182 down to 86..183)
# and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("85:85").Insert()

$ws.Range("A85").Value = 10
$ws.Range("B85").Value = "Vega Modelo de Temuco"
$ws.Range("C85").Value = "La Araucanía"
$ws.Range("D85").Value = 44580
$ws.Range("E85").Value = 9
$ws.Range("F85").Value = 100112052
$ws.Range("G85").Value = "Albahaca"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 70
$ws.Range("K85").Value = 4000
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = 4429
$ws.Range("N85").Value = "$/paquete"
$ws.Range("O85").Value = "Región del Maule"
$ws.Range("P85").Value = 4429
$ws.Range("Q85").Value = 1
$ws.Range("R85").Value = "Hortaliza"
